# "Generate Report for Handoff"
#
# A new handoff-xliff batch was generated, refreshing the handoff
# timestamps for the set of files that were re-handed-off (rows 7, 9,
# 10, 11, 13, 14 -- i.e. every "Ready for handoff" row except the two
# that are mid-cycle already, rows 8 and 12):
#
#   - "Latest HO Xliff Generate Date" on the Overview sheet, and the
#     de-de sheet's "Latest Handoff Datetime", shared the same source
#     timestamp (2016-08-16 00:19:42) and both move to
#     2016-08-16 00:20:15.
#   - the zh-cn sheet's "Latest Handoff Datetime" had its own timestamp
#     (2016-08-16 00:19:37) and moves to 2016-08-16 00:20:05.
#   - those same rows now carry the "ht" (handoff type) marker in the
#     Priority column, on both language sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 13, 14)

# Overview sheet: bump the latest HO xliff generate date for the files
# that were just handed off again.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-16 00:20:15"
}

# de-de sheet shared the very same "generate date" string as Overview.
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-08-16 00:20:15"
    $dede.Range("E$r").Value = "ht"
}

# zh-cn sheet had its own distinct timestamp for this batch.
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-08-16 00:20:05"
    $zhcn.Range("E$r").Value = "ht"
}
